$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the top of the data block (rows 518-519),
# pushing the existing rows 518-557 down to 520-559.
$ws.Range("A518:A519").EntireRow.Insert()

# New row 518: "Primera" quality entry for the latest week (date 45265)
$ws.Cells.Item(518, 1).Value = 1
$ws.Cells.Item(518, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(518, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(518, 4).Value = 45265
$ws.Cells.Item(518, 5).Value = 15
$ws.Cells.Item(518, 6).Value = 100112043
$ws.Cells.Item(518, 7).Value = "Pepino ensalada"
$ws.Cells.Item(518, 8).Value = "Sin especificar"
$ws.Cells.Item(518, 9).Value = "Primera"
$ws.Cells.Item(518, 10).Value = 120
$ws.Cells.Item(518, 11).Value = 11000
$ws.Cells.Item(518, 12).Value = 12000
$ws.Cells.Item(518, 13).Value = 11500
$ws.Cells.Item(518, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(518, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(518, 16).Value = 164
$ws.Cells.Item(518, 17).Value = 70
$ws.Cells.Item(518, 18).Value = "Hortaliza"

# New row 519: "Segunda" quality entry for the latest week (date 45265)
$ws.Cells.Item(519, 1).Value = 1
$ws.Cells.Item(519, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(519, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(519, 4).Value = 45265
$ws.Cells.Item(519, 5).Value = 15
$ws.Cells.Item(519, 6).Value = 100112043
$ws.Cells.Item(519, 7).Value = "Pepino ensalada"
$ws.Cells.Item(519, 8).Value = "Sin especificar"
$ws.Cells.Item(519, 9).Value = "Segunda"
$ws.Cells.Item(519, 10).Value = 140
$ws.Cells.Item(519, 11).Value = 8000
$ws.Cells.Item(519, 12).Value = 9000
$ws.Cells.Item(519, 13).Value = 8500
$ws.Cells.Item(519, 14).Value = "`$/caja 100 unidades"
$ws.Cells.Item(519, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(519, 16).Value = 85
$ws.Cells.Item(519, 17).Value = 100
$ws.Cells.Item(519, 18).Value = "Hortaliza"
